# Generate Report for Handoff
# Adds two new localization entries (197fab40... and dcdb0ae3...) to the
# Overview / zh-cn / de-de sheets of the localization-status workbook,
# inserting "197fab40" before the existing "571c7103" row and appending
# "dcdb0ae3" after it, on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Remove existing hyperlinks so they can be rewritten cleanly in the
# correct final order (this engine's Hyperlinks collection is append-only
# per-ref, so the safest path is: wipe, then re-add everything in order).
$ws1.Hyperlinks.Delete()

# Row 2 (unchanged data, rewritten so the hyperlink can be re-added)
$ws1.Range("A2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-21 22:35:35"
$ws1.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3 (NEW: 197fab40)
$ws1.Range("A3").Value = "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-21 22:36:37"
$ws1.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 4 (was row 3: 571c7103)
$ws1.Range("A4").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-03-21 22:34:46"
$ws1.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 5 (NEW: dcdb0ae3)
$ws1.Range("A5").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-03-21 22:36:37"
$ws1.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Re-create the File Name hyperlinks in final row order
$ws1.Range("A2:A5").Style = "HyperLink"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7d9cae1f907f2f0e19d72d98df700e24046a7a28/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ca0cc54eaaca092c89435a4ef7fd610b26fd008/e2e/197fab40-7705-4dd6-a3a2-ec57183f1ea7.md", "", "", "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/571c7103-8660-4f11-8e8b-df8803d0e27d.md", "", "", "571c7103-8660-4f11-8e8b-df8803d0e27d.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/652554d13b5908c2af336ef3206db248a3d31e97/e2e/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md", "", "", "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md") | Out-Null

Write-Host "Sheet1 (Overview) done"

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

# Row 2 (unchanged data, rewritten so the hyperlinks can be re-added)
$ws2.Range("A2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-21 22:35:31"
$ws2.Range("F2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.md"
$ws2.Range("G2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-21 22:35:48"
$ws2.Range("J2").Value = "Include"

# Row 3 (NEW: 197fab40)
$ws2.Range("A3").Value = "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-21 22:36:31"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("J3").Value = "Include"

# Row 4 (was row 3: 571c7103)
$ws2.Range("A4").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-21 22:34:43"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("J4").Value = "Include"

# Row 5 (NEW: dcdb0ae3)
$ws2.Range("A5").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.zh-cn.xlf"
$ws2.Range("E5").Value = "2016-03-21 22:36:31"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("J5").Value = "Include"

# Datetime-style columns E/H keep the workbook's datetime display format
$ws2.Range("E2:E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H2:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlink-style columns (A, D, F, G)
$ws2.Range("A2:A5").Style = "HyperLink"
$ws2.Range("D2:D5").Style = "HyperLink"
$ws2.Range("F2").Style = "HyperLink"
$ws2.Range("G2").Style = "HyperLink"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7d9cae1f907f2f0e19d72d98df700e24046a7a28/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8da87dd7580ef39e6e092d7e540a6cab6d528832/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a734c2cb4d4d89f1fdb4589337d30735c022bd58/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/01f75f6f65a8b3d6abb033821e7cfdf48a1b05a6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0ca0cc54eaaca092c89435a4ef7fd610b26fd008/e2e/197fab40-7705-4dd6-a3a2-ec57183f1ea7.md", "", "", "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3cdd7ae476fb505d3659423fc70cc6f9e768530a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.zh-cn.xlf", "", "", "197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/571c7103-8660-4f11-8e8b-df8803d0e27d.md", "", "", "571c7103-8660-4f11-8e8b-df8803d0e27d.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86c5f1c3fcf4eeb680ca7d9a65b3ee542c89a9bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.zh-cn.xlf", "", "", "571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/652554d13b5908c2af336ef3206db248a3d31e97/e2e/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md", "", "", "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2a96675ae8f83a8494f956f11785ee6209dc9cf2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.zh-cn.xlf", "", "", "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.zh-cn.xlf") | Out-Null

Write-Host "Sheet2 (zh-cn) done"

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

# Row 2 (unchanged data, rewritten so the hyperlinks can be re-added)
$ws3.Range("A2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-21 22:35:35"
$ws3.Range("F2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.md"
$ws3.Range("G2").Value = "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-21 22:35:54"
$ws3.Range("J2").Value = "Include"

# Row 3 (NEW: 197fab40)
$ws3.Range("A3").Value = "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-21 22:36:37"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("J3").Value = "Include"

# Row 4 (was row 3: 571c7103)
$ws3.Range("A4").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-21 22:34:46"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("J4").Value = "Include"

# Row 5 (NEW: dcdb0ae3)
$ws3.Range("A5").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.de-de.xlf"
$ws3.Range("E5").Value = "2016-03-21 22:36:37"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("J5").Value = "Include"

# Datetime-style columns E/H keep the workbook's datetime display format
$ws3.Range("E2:E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H2:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlink-style columns (A, D, F, G)
$ws3.Range("A2:A5").Style = "HyperLink"
$ws3.Range("D2:D5").Style = "HyperLink"
$ws3.Range("F2").Style = "HyperLink"
$ws3.Range("G2").Style = "HyperLink"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7d9cae1f907f2f0e19d72d98df700e24046a7a28/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86ebc92335ca81ee9d52dae18a8f7f0bd62a2f34/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/babad688c7f84e7d63a0266026211f6fdafe15bb/e2e/64d947ee-a767-493c-8a90-ea0403e5866b.md", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e7faf877098ed28ec94e1703d5ffd504aa2d4761/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf", "", "", "64d947ee-a767-493c-8a90-ea0403e5866b.bd75e25f15b28351a1a39d6513bbb031ca67d6f8.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cd10e49cff1ba08cedbc94209bc201351e977369/e2e/197fab40-7705-4dd6-a3a2-ec57183f1ea7.md", "", "", "197fab40-7705-4dd6-a3a2-ec57183f1ea7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd10e49cff1ba08cedbc94209bc201351e977369/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.de-de.xlf", "", "", "197fab40-7705-4dd6-a3a2-ec57183f1ea7.ff9333eca2a7c0f89f719f5b6f586b21d44f70c4.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/571c7103-8660-4f11-8e8b-df8803d0e27d.md", "", "", "571c7103-8660-4f11-8e8b-df8803d0e27d.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34ae5c474cf2361ac996412ee2a82e4e64ab8941/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.de-de.xlf", "", "", "571c7103-8660-4f11-8e8b-df8803d0e27d.f54e4f0a38172b86200b0bb35895fd010cca9747.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/70d71f289588e95054e838f0513a57fc3eb57d68/e2e/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md", "", "", "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70d71f289588e95054e838f0513a57fc3eb57d68/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.de-de.xlf", "", "", "dcdb0ae3-dd2f-4bcf-aad6-b54db1786366.fc7c16d2e7cf0b5634a92a7d8b94b4eb3218ae84.de-de.xlf") | Out-Null

Write-Host "Sheet3 (de-de) done"
